$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the swap between row 23 and row 24
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr23 = "$col" + "23"
    $addr24 = "$col" + "24"

    $val23 = $ws.Range($addr23).Value()
    $val24 = $ws.Range($addr24).Value()

    $ws.Range($addr23).Value = $val24
    $ws.Range($addr24).Value = $val23
}
